# Weekly refresh: a new week's price observation is inserted at the top of
# the data block (row 88), pushing every existing record down by one row.
# The previously-last record (old row 206) becomes the new last record
# (row 207), growing the used range from A1:T206 to A1:T207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 88; Excel shifts rows 88-206 down to
# 89-207 and carries their formatting (e.g. the date style on column D).
$ws.Rows.Item(88).Insert()

# Populate the newly blank row 88 with this week's observation.
$ws.Cells.Item(88, 1).Value  = 4
$ws.Cells.Item(88, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(88, 3).Value  = "Los Lagos"
$ws.Cells.Item(88, 4).Value  = 44579
$ws.Cells.Item(88, 5).Value  = 10
$ws.Cells.Item(88, 6).Value  = "Fruta"
$ws.Cells.Item(88, 7).Value  = 100104
$ws.Cells.Item(88, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(88, 9).Value  = 100104005
$ws.Cells.Item(88, 10).Value = "Pera"
$ws.Cells.Item(88, 11).Value = "Packham's Triumph"
$ws.Cells.Item(88, 12).Value = "Primera"
$ws.Cells.Item(88, 13).Value = 500
$ws.Cells.Item(88, 14).Value = 14000
$ws.Cells.Item(88, 15).Value = 15000
$ws.Cells.Item(88, 16).Value = 14500
$ws.Cells.Item(88, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(88, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(88, 19).Value = 967
$ws.Cells.Item(88, 20).Value = 15
